$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# setup some test row
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(5, $c).Value = $c * 100
}
Write-Host "Before delete:"
for ($c = 1; $c -le 6; $c++) {
    Write-Host "  col $c :" $ws.Cells.Item(5,$c).Value()
}

# delete A5 with shift left
$rng = $ws.Range("A5")
$rng.Delete(-4159)  # xlShiftToLeft = -4159
Write-Host "After delete shiftleft:"
for ($c = 1; $c -le 6; $c++) {
    Write-Host "  col $c :" $ws.Cells.Item(5,$c).Value()
}
